# Fill in rows 14 and 15 of the "Spieltabelle" sheet (poker session tracking
# log) with the data for sessions 12 and 13, which were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spieltabelle")

# --- Row 14 (Index 12) ---------------------------------------------------
$ws.Range("B14").Value = "Cashgame"
$ws.Range("C14").Value = "sc.ch"
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 0.8
$ws.Range("F14").Value = 0.01
$ws.Range("G14").Formula = "=F14*2"
$ws.Range("H14").Value = 45966
$ws.Range("I14").Value = 0.87152777777777779
$ws.Range("J14").Value = 45966
$ws.Range("K14").Value = 0.88888888888888884
$ws.Range("L14").Formula = "=(J14+K14)-(H14+I14)"
$ws.Range("M14").Value = 1.5
$ws.Range("N14").Value = 0
$ws.Range("O14").Formula = "=N14-M14"
$ws.Range("P14").Formula = '=LET(out, SUBSTITUTE(O14,"CHF ",""), IF(out="-", "", IFERROR( NUMBERVALUE(out) / (L14*24), "" )))'
$ws.Range("Q14").Value = 1.5
$ws.Range("R14").Value = 3
$ws.Range("S14").Value = 0.38
$ws.Range("T14").Value = "Nichts"

# --- Row 15 (Index 13) ---------------------------------------------------
$ws.Range("B15").Value = "Cashgame"
$ws.Range("C15").Value = "sc.ch"
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 0.8
$ws.Range("F15").Value = 0.01
$ws.Range("G15").Formula = "=F15*2"
$ws.Range("H15").Value = 45966
$ws.Range("I15").Value = 0.88888888888888884
$ws.Range("J15").Value = 45966
$ws.Range("K15").Value = 0.89097222222222228
$ws.Range("L15").Formula = "=(J15+K15)-(H15+I15)"
$ws.Range("M15").Value = 1.5
$ws.Range("N15").Value = 0
$ws.Range("O15").Formula = "=N15-M15"
$ws.Range("P15").Formula = '=LET(out, SUBSTITUTE(O15,"CHF ",""), IF(out="-", "", IFERROR( NUMBERVALUE(out) / (L15*24), "" )))'
$ws.Range("Q15").Value = 1.5
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 0.15
$ws.Range("T15").Value = "Nichts"

# Remark/opponent notes - set in this specific order so that new shared
# strings are appended to the shared-string table in the same sequence as
# in the target workbook ("Keine brauchbaren Karten", "Ein ASo verlor gegen
# AA", "Keine Aussage möglich").
$ws.Range("U14").Value = "Keine brauchbaren Karten"
$ws.Range("U15").Value = "Ein ASo verlor gegen AA"
$ws.Range("V14").Value = "Keine Aussage möglich"
$ws.Range("V15").Value = "Keine Aussage möglich"

# Recalculate formulas now that their inputs are populated.
$wb.Application.CalculateFull()

# Update the selection/scroll position to match where the editor ended up.
$ws.Range("U14").Select()
